$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.865.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -7.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.708.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5168"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -12.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.92%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06213"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -8.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07350"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.699.57"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -7.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.489"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5812"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.937.39"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.23%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -11.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.71"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -12.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.904.77"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.019"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.68"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "186.60"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -10.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.273"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.69"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.583"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1154"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -9.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.23"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05888"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.348"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.472"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.442"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.663"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9898"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6024"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.408"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.688"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.103.74"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01600"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8672"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.908"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -8.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.67"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.844.23"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.23"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000105"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4380"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05251"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.913"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.03%  "
